# Update "Datos actualizados" timestamp (row 1) and refresh the
# provincias/ciudades COVID table: update case counts for the cities whose
# figures changed, then re-sort the whole table (rows 4-61) descending by
# "Casos totales" (column B), matching the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "last updated" timestamp.
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 14:16"

# Final table contents (already sorted descending by Casos totales), one
# entry per data row: Ciudad, Casos totales, Casos activos, Recuperados, Muertes.
$data = @(
  @('Madrid', 7165, 1186, 5351, 628),
  @('Cataluña', 3270, 3, 3185, 82),
  @('Araba/Alava', 703, 21, 655, 48),
  @('Valencia/Valencia', 627, 12, 600, 15),
  @('Navarra', 593, 2, 585, 6),
  @('Bizkaia/Vizcaya', 539, 21, 522, 17),
  @('La Rioja', 497, 13, 477, 7),
  @('Ciudad Real', 400, 8, 364, 28),
  @('Alacant/Alicante', 372, 7, 348, 17),
  @('Toledo', 370, 15, 336, 19),
  @('Malaga', 361, 0, 348, 13),
  @('Asturias', 344, 12, 329, 3),
  @('Albacete', 327, 8, 291, 28),
  @('A Coruña', 270, 0, 267, 3),
  @('Guadalajara', 237, 2, 231, 4),
  @('Zaragoza', 224, 0, 210, 14),
  @('Gipuzkoa/Guipuzcoa', 223, 21, 217, 6),
  @('Murcia', 215, 1, 214, 0),
  @('Burgos', 208, 21, 175, 12),
  @('Caceres', 206, 2, 194, 10),
  @('Illes Balears', 203, 10, 189, 4),
  @('Salamanca', 203, 10, 180, 13),
  @('Valladolid', 201, 3, 193, 5),
  @('Pontevedra', 193, 0, 191, 2),
  @('Tenerife', 192, 4, 185, 3),
  @('Granada', 176, 0, 169, 7),
  @('Aragon', 174, 0, 163, 11),
  @('Illes Balears*', 169, 6, 161, 2),
  @('Leon', 164, 3, 156, 5),
  @('Segovia', 156, 6, 140, 10),
  @('Cantabria', 144, 10, 133, 1),
  @('Sevilla', 133, 1, 131, 1),
  @('Castello/Castellon', 104, 1, 102, 1),
  @('Cordoba', 101, 0, 101, 0),
  @('Cuenca', 94, 5, 84, 5),
  @('Badajoz', 91, 5, 86, 0),
  @('Jaen', 87, 0, 85, 2),
  @('Cadiz', 84, 0, 84, 0),
  @('Soria', 80, 4, 71, 5),
  @('Gran Canaria', 70, 0, 69, 1),
  @('Ourense', 63, 0, 63, 0),
  @('Avila', 63, 4, 55, 4),
  @('Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
  @('Lugo', 46, 0, 45, 1),
  @('Zamora', 46, 3, 42, 1),
  @('Almeria', 37, 0, 37, 0),
  @('Teruel', 27, 0, 26, 1),
  @('Palencia', 26, 1, 25, 0),
  @('Melilla', 24, 0, 24, 0),
  @('Huesca', 24, 0, 24, 0),
  @('Huelva', 23, 0, 23, 0),
  @('Fuerteventura', 12, 0, 12, 0),
  @('La Palma', 7, 0, 7, 0),
  @('Arroyo de la Luz', 7, 0, 7, 0),
  @('Ceuta', 5, 0, 5, 0),
  @('Lanzarote', 3, 0, 3, 0),
  @('La Gomera', 3, 2, 1, 0),
  @('El Hierro', 1, 0, 1, 0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
